$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9817641377449036
$ws.Range("B1").Value = 2.201106071472168
$ws.Range("C1").Value = 5.091033458709717
$ws.Range("D1").Value = 1.843563079833984
$ws.Range("E1").Value = 1.300118803977966
